# Add two new trailing columns, I ("I0") and J ("IF"), to the sheet.
# The sheet's used range grows from A1:H64 to A1:J64.
#
# Header cells (I1/J1) need the same style as the other header cells
# (bold font, thin border, centered/top alignment - style index 1 in
# the original workbook). Copy/PasteSpecial(formats) from the existing
# H1 header reuses that exact style instead of minting a near-duplicate
# one, matching what Excel itself would do when you fill a header row
# across by copy-paste.
#
# Data cells (I2:J64/I64:J64) get plain numeric values with no special
# style, matching columns C:H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers ---------------------------------------------------------
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data (rows 2..64) ------------------------------------------------
$I0 = @(8,8,8,9,8,9,8,9,9,9,9,10,9,9,10,9,9,7,9,8,9,9,9,9,10,9,9,9,10,9,9,9,9,8,8,9,9,9,9,8,9,8,9,9,9,9,9,9,9,9,8,9,8,9,8,8,9,9,6,6,5,3,3)
$IF = @(9,8,8,9,9,9,8,9,9,9,9,10,9,9,10,9,9,7,9,8,9,9,9,9,10,9,9,9,10,9,9,9,9,8,8,9,9,9,10,8,9,8,9,9,9,9,9,9,9,9,8,9,8,9,8,8,9,9,6,6,5,3,3)

for ($i = 0; $i -lt $I0.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $I0[$i]
    $ws.Cells.Item($row, 10).Value = $IF[$i]
}
